$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.218.57'
$ws.Range("E2").Value = '  +0.74%  '

# Row 3
$ws.Range("D3").Value = '1.804.43'
$ws.Range("E3").Value = '  +2.49%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = "'338.92"
$ws.Range("E5").Value = '  +0.54%  '

# Row 6
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = '  +0.06%  '

# Row 7
$ws.Range("D7").Value = "'0.4917"
$ws.Range("E7").Value = '  +29.95%  '

# Row 8
$ws.Range("D8").Value = "'0.3731"
$ws.Range("E8").Value = '  +11.15%  '

# Row 9
$ws.Range("D9").Value = "'45.60"
$ws.Range("E9").Value = '  -0.05%  '

# Row 10
$ws.Range("D10").Value = "'0.07761"
$ws.Range("E10").Value = '  +7.84%  '

# Row 11
$ws.Range("E11").Value = '  +2.54%  '

# Row 12
$ws.Range("D12").Value = "'22.67"
$ws.Range("E12").Value = '  +1.34%  '

# Row 13
$ws.Range("D13").Value = "'1.002"
$ws.Range("E13").Value = '  +0.01%  '

# Row 14
$ws.Range("D14").Value = "'6.337"
$ws.Range("E14").Value = '  +2.30%  '

# Row 15
$ws.Range("D15").Value = "'7.339"
$ws.Range("E15").Value = '  +1.88%  '

# Row 16
$ws.Range("D16").Value = '1.794.28'
$ws.Range("E16").Value = '  +1.98%  '

# Row 17
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = '  +4.31%  '

# Row 18
$ws.Range("D18").Value = "'0.06743"
$ws.Range("E18").Value = '  +2.49%  '

# Row 19
$ws.Range("D19").Value = "'82.42"
$ws.Range("E19").Value = '  +2.49%  '

# Row 20
$ws.Range("D20").Value = "'0.9993"
$ws.Range("E20").Value = '  -0.01%  '

# Row 21
$ws.Range("E21").Value = '  +2.70%  '

# Row 22
$ws.Range("D22").Value = "'6.431"
$ws.Range("E22").Value = '  +2.37%  '

# Row 23
$ws.Range("D23").Value = '28.194.77'
$ws.Range("E23").Value = '  +0.68%  '

# Row 24
$ws.Range("D24").Value = "'12.02"
$ws.Range("E24").Value = '  +2.67%  '

# Row 25
$ws.Range("D25").Value = "'2.399"
$ws.Range("E25").Value = '  +1.13%  '

# Row 26
$ws.Range("D26").Value = "'20.92"
$ws.Range("E26").Value = '  +5.23%  '

# Row 27
$ws.Range("D27").Value = "'2.420"
$ws.Range("E27").Value = '  +3.49%  '

# Row 28
$ws.Range("D28").Value = "'151.43"
$ws.Range("E28").Value = '  -1.21%  '

# Row 29
$ws.Range("D29").Value = '2.004.48'
$ws.Range("E29").Value = '  +2.25%  '

# Row 30
$ws.Range("D30").Value = "'134.63"
$ws.Range("E30").Value = '  +2.09%  '

# Row 31
$ws.Range("D31").Value = "'1.277"
$ws.Range("E31").Value = '  +1.58%  '

# Row 32
$ws.Range("D32").Value = "'4.044"
$ws.Range("E32").Value = '  +0.68%  '

# Row 33
$ws.Range("D33").Value = "'0.09921"
$ws.Range("E33").Value = '  +12.57%  '

# Row 34
$ws.Range("D34").Value = "'5.968"
$ws.Range("E34").Value = '  +2.78%  '

# Row 35
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = "'12.26"
$ws.Range("E35").Value = '  +0.09%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = "'0.02389"
$ws.Range("E36").Value = '  +2.00%  '

# Row 37
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.06389"
$ws.Range("E37").Value = '  +3.33%  '

# Row 38
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.2225"
$ws.Range("E38").Value = '  +5.05%  '

# Row 39
$ws.Range("D39").Value = "'0.6713"
$ws.Range("E39").Value = '  +1.22%  '

# Row 40
$ws.Range("D40").Value = "'5.257"
$ws.Range("E40").Value = '  +1.77%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'1.224"
$ws.Range("E41").Value = '  +1.00%  '

# Row 42
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = "'1.486"
$ws.Range("E42").Value = '  +2.42%  '

# Row 43
$ws.Range("D43").Value = "'8.127"
$ws.Range("E43").Value = '  +1.03%  '

# Row 44
$ws.Range("D44").Value = "'14.16"
$ws.Range("E44").Value = '  +2.18%  '

# Row 45
$ws.Range("D45").Value = "'0.9991"

# Row 46
$ws.Range("D46").Value = "'0.6183"
$ws.Range("E46").Value = '  +1.87%  '

# Row 47
$ws.Range("E47").Value = '  +1.23%  '

# Row 48
$ws.Range("D48").Value = "'129.36"
$ws.Range("E48").Value = '  -0.21%  '

# Row 49
$ws.Range("E49").Value = '  +2.08%  '

# Row 50
$ws.Range("D50").Value = "'1.178"
$ws.Range("E50").Value = '  -0.79%  '

# Row 51
$ws.Range("D51").Value = "'0.07121"
$ws.Range("E51").Value = '  -0.58%  '
